# Auto-generated edit script: refresh scraped horarios data (Linea 141, 10/01/2026)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 19:47:58"
$ws.Cells.Item(3,1).Value = "Total filas: 121"
$data = @(
  @("16:46:42","16:47","15_ABASTO",1,"LP1912"),
  @("16:50:41","16:50","10_OLMOS",0,"LP1912"),
  @("16:50:41","16:51","15_ABASTO",1,"LP1912"),
  @("16:52:37","16:53","15_ABASTO",1,"LP1912"),
  @("16:46:42","16:53","10_OLMOS",7,"LP1912"),
  @("16:46:42","16:56","215C_EL PATO",10,"LP1912"),
  @("16:46:42","17:01","16_SANTA ANA",15,"LP1912"),
  @("16:46:42","17:03","23_HERNANDEZ",17,"LP1912"),
  @("16:46:42","17:04","14_ABASTO",18,"LP1912"),
  @("16:46:42","17:07","15_ABASTO",21,"LP1912"),
  @("16:46:42","17:13","23_HERNANDEZ",27,"LP1912"),
  @("16:46:42","17:14","10_OLMOS",28,"LP1912"),
  @("17:13:30","17:16","10_OLMOS",3,"LP1912"),
  @("16:46:42","17:17","17_ROMERO",31,"LP1912"),
  @("16:50:41","17:17","23_HERNANDEZ",27,"LP1912"),
  @("16:52:37","17:20","23_HERNANDEZ",28,"LP1912"),
  @("16:46:42","17:23","16_SANTA ANA",37,"LP1912"),
  @("16:46:42","17:24","11_ETCHEVERRY",38,"LP1912"),
  @("17:13:30","17:27","15_ABASTO",14,"LP1912"),
  @("17:13:30","17:33","23_HERNANDEZ",20,"LP1912"),
  @("16:50:41","17:34","10_OLMOS",44,"LP1912"),
  @("16:46:42","17:35","16_P MOR-SANTA ANA",49,"LP1912"),
  @("16:52:37","17:36","27_EL RETIRO",44,"LP1912"),
  @("17:13:30","17:37","27_EL RETIRO",24,"LP1912"),
  @("16:46:42","17:38","17X38_ROMERO",52,"LP1912"),
  @("17:35:09","17:39","27_EL RETIRO",4,"LP1912"),
  @("17:13:30","17:41","23_HERNANDEZ",28,"LP1912"),
  @("16:46:42","17:44","215B_EL PATO",58,"LP1912"),
  @("17:35:09","17:45","215B_EL PATO",10,"LP1912"),
  @("16:50:41","17:47","16_SANTA ANA",57,"LP1912"),
  @("16:46:42","17:48","27_EL RETIRO",62,"LP1912"),
  @("16:50:41","17:49","27_EL RETIRO",59,"LP1912"),
  @("16:46:42","17:50","215_EL PELIGRO",64,"LP1912"),
  @("17:47:22","17:51","215B_EL PATO",4,"LP1912"),
  @("16:52:37","17:51","215_EL PELIGRO",59,"LP1912"),
  @("17:47:22","17:54","10_OLMOS",7,"LP1912"),
  @("17:35:09","18:00","16_SANTA ANA",25,"LP1912"),
  @("16:46:42","18:02","17_ROMERO",76,"LP1912"),
  @("16:52:37","18:03","17_ROMERO",71,"LP1912"),
  @("17:35:09","18:03","23_HERNANDEZ",28,"LP1912"),
  @("16:46:42","18:04","14_ABASTO",78,"LP1912"),
  @("17:35:09","18:05","14_ABASTO",30,"LP1912"),
  @("18:10:41","18:11","10_OLMOS",1,"LP1912"),
  @("18:10:41","18:11","16_SANTA ANA",1,"LP1912"),
  @("16:52:37","18:14","10_OLMOS",82,"LP1912"),
  @("17:47:22","18:21","16_SANTA ANA",34,"LP1912"),
  @("16:46:42","18:24","11_ETCHEVERRY",98,"LP1912"),
  @("17:35:09","18:25","11_ETCHEVERRY",50,"LP1912"),
  @("17:13:30","18:27","15_ABASTO",74,"LP1912"),
  @("17:47:22","18:31","16_SANTA ANA",44,"LP1912"),
  @("17:35:09","18:31","23_HERNANDEZ",56,"LP1912"),
  @("17:54:43","18:33","23_HERNANDEZ",39,"LP1912"),
  @("16:46:42","18:34","14X44_ABASTO",108,"LP1912"),
  @("16:46:42","18:38","17X38_ROMERO",112,"LP1912"),
  @("16:46:42","18:41","16_P MOR-SANTA ANA",115,"LP1912"),
  @("17:13:30","18:41","14_ABASTO",88,"LP1912"),
  @("17:47:22","18:44","14_ABASTO",57,"LP1912"),
  @("17:35:09","18:45","14_ABASTO",70,"LP1912"),
  @("18:44:34","18:47","14_ABASTO",3,"LP1912"),
  @("17:35:09","18:51","15_ABASTO",76,"LP1912"),
  @("18:52:04","18:52","15_ABASTO",0,"LP1912"),
  @("17:54:43","18:53","16_SANTA ANA",59,"LP1912"),
  @("17:35:09","18:59","10_OLMOS",84,"LP1912"),
  @("17:13:30","19:01","17_ROMERO",108,"LP1912"),
  @("18:10:41","19:03","23_HERNANDEZ",53,"LP1912"),
  @("18:52:04","19:04","23_HERNANDEZ",12,"LP1912"),
  @("18:31:18","19:05","16_SANTA ANA",34,"LP1912"),
  @("17:13:30","19:11","81_EL PELIGRO",118,"LP1912"),
  @("18:10:41","19:14","14_ABASTO",64,"LP1912"),
  @("18:52:04","19:15","14_ABASTO",23,"LP1912"),
  @("18:44:34","19:17","16_SANTA ANA",33,"LP1912"),
  @("17:47:22","19:17","27_EL RETIRO",90,"LP1912"),
  @("17:35:09","19:19","27_EL RETIRO",104,"LP1912"),
  @("17:54:43","19:20","215C_EL PATO",86,"LP1912"),
  @("17:35:09","19:21","215C_EL PATO",106,"LP1912"),
  @("19:11:59","19:22","27_EL RETIRO",11,"LP1912"),
  @("18:44:34","19:26","27_EL RETIRO",42,"LP1912"),
  @("19:11:59","19:28","225_GOMEZ",17,"LP1912"),
  @("17:35:09","19:29","225_GOMEZ",114,"LP1912"),
  @("17:54:43","19:30","27_EL RETIRO",96,"LP1912"),
  @("17:54:43","19:30","215_EL PELIGRO",96,"LP1912"),
  @("17:35:09","19:31","215_EL PELIGRO",116,"LP1912"),
  @("18:10:41","19:31","27_EL RETIRO",81,"LP1912"),
  @("18:44:34","19:33","23_HERNANDEZ",49,"LP1912"),
  @("18:31:18","19:34","23_HERNANDEZ",63,"LP1912"),
  @("19:11:59","19:38","16_SANTA ANA",27,"LP1912"),
  @("19:11:59","19:39","17X38_ROMERO",28,"LP1912"),
  @("19:11:59","19:40","16_SANTA ANA",29,"LP1912"),
  @("17:47:22","19:40","17X38_ROMERO",113,"LP1912"),
  @("19:35:31","19:41","16_SANTA ANA",6,"LP1912"),
  @("19:11:59","19:43","11_ETCHEVERRY",32,"LP1912"),
  @("17:47:22","19:44","11_ETCHEVERRY",117,"LP1912"),
  @("18:31:18","19:46","11_ETCHEVERRY",75,"LP1912"),
  @("19:47:58","19:47","16_SANTA ANA",0,"LP1912"),
  @("19:11:59","19:50","81_EL PELIGRO",39,"LP1912"),
  @("17:54:43","19:51","81_EL PELIGRO",117,"LP1912"),
  @("18:10:41","19:58","14X44_ABASTO",108,"LP1912"),
  @("18:31:18","19:59","14X44_ABASTO",88,"LP1912"),
  @("18:10:41","20:00","215C_EL PATO",110,"LP1912"),
  @("19:47:58","20:00","16_SANTA ANA",13,"LP1912"),
  @("18:31:18","20:01","215C_EL PATO",90,"LP1912"),
  @("19:47:58","20:02","17X38_ROMERO",15,"LP1912"),
  @("19:11:59","20:04","23_HERNANDEZ",53,"LP1912"),
  @("19:47:58","20:09","23_HERNANDEZ",22,"LP1912"),
  @("19:35:31","20:10","23_HERNANDEZ",35,"LP1912"),
  @("19:11:59","20:13","11_ETCHEVERRY",62,"LP1912"),
  @("18:31:18","20:14","11_ETCHEVERRY",103,"LP1912"),
  @("19:11:59","20:25","15_ABASTO",74,"LP1912"),
  @("18:31:18","20:26","15_ABASTO",115,"LP1912"),
  @("18:44:34","20:28","10_OLMOS",104,"LP1912"),
  @("18:31:18","20:29","10_OLMOS",118,"LP1912"),
  @("19:11:59","20:43","215B_EL PATO",92,"LP1912"),
  @("19:11:59","20:44","17X38_ROMERO",93,"LP1912"),
  @("18:52:04","20:44","215B_EL PATO",112,"LP1912"),
  @("18:52:04","20:45","17X38_ROMERO",113,"LP1912"),
  @("19:11:59","20:52","23_HERNANDEZ",77,"LP1912"),
  @("19:11:59","21:01","215A_EL PATO",110,"LP1912"),
  @("19:11:59","21:02","27_EL RETIRO",111,"LP1912"),
  @("19:47:58","21:06","27_EL RETIRO",79,"LP1912"),
  @("19:35:31","21:10","27_EL RETIRO",95,"LP1912"),
  @("19:35:31","21:23","10_OLMOS",108,"LP1912")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $r = $r + 1
}

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 19:47:58"
$ws.Cells.Item(3,1).Value = "Total filas: 15"
$data = @(
  @("16:46:42","16:56","215C_EL PATO",10,"LP1912"),
  @("16:46:42","17:44","215B_EL PATO",58,"LP1912"),
  @("17:35:09","17:45","215B_EL PATO",10,"LP1912"),
  @("16:46:42","17:50","215_EL PELIGRO",64,"LP1912"),
  @("16:52:37","17:51","215_EL PELIGRO",59,"LP1912"),
  @("17:47:22","17:51","215B_EL PATO",4,"LP1912"),
  @("17:54:43","19:20","215C_EL PATO",86,"LP1912"),
  @("17:35:09","19:21","215C_EL PATO",106,"LP1912"),
  @("17:54:43","19:30","215_EL PELIGRO",96,"LP1912"),
  @("17:35:09","19:31","215_EL PELIGRO",116,"LP1912"),
  @("18:10:41","20:00","215C_EL PATO",110,"LP1912"),
  @("18:31:18","20:01","215C_EL PATO",90,"LP1912"),
  @("19:11:59","20:43","215B_EL PATO",92,"LP1912"),
  @("18:52:04","20:44","215B_EL PATO",112,"LP1912"),
  @("19:11:59","21:01","215A_EL PATO",110,"LP1912")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $r = $r + 1
}

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 19:47:58"
$ws.Cells.Item(3,1).Value = "Total filas: 14"
$data = @(
  @("16:52:37","16:57","215C_LA PLATA",5,"L6203"),
  @("16:50:41","16:58","215C_LA PLATA",8,"L6203"),
  @("16:46:42","16:59","215C_LA PLATA",13,"L6203"),
  @("16:46:42","18:21","215C_LA PLATA",95,"L6203"),
  @("17:13:30","18:22","215C_LA PLATA",69,"L6203"),
  @("17:13:30","19:11","215B_LP-P MOR-1 Y 57",118,"L6173"),
  @("17:47:22","19:12","215B_LP-P MOR-1 Y 57",85,"L6173"),
  @("18:31:18","19:13","215B_LP-P MOR-1 Y 57",42,"L6173"),
  @("18:52:04","19:14","215B_LP-P MOR-1 Y 57",22,"L6173"),
  @("18:10:41","19:15","215B_LP-P MOR-1 Y 57",65,"L6173"),
  @("18:44:34","20:41","215C_LA PLATA",117,"L6203"),
  @("18:52:04","20:42","215C_LA PLATA",110,"L6203"),
  @("19:47:58","21:33","215C_LA PLATA",106,"L6203"),
  @("19:35:31","21:34","215C_LA PLATA",119,"L6203")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $r = $r + 1
}
